# Insert two new data rows (973, 974) into the "Pera" dataset on the
# active sheet. This shifts all existing rows 973:1049 down to 975:1051,
# growing the used range from A1:T1049 to A1:T1051.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("973:974").Insert()

# --- New row 973 -----------------------------------------------------
$ws.Cells.Item(973, 1).Value  = 9
$ws.Cells.Item(973, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(973, 3).Value  = "Metropolitana"
$ws.Cells.Item(973, 4).Value  = 45013
$ws.Cells.Item(973, 5).Value  = 13
$ws.Cells.Item(973, 6).Value  = "Fruta"
$ws.Cells.Item(973, 7).Value  = 100104
$ws.Cells.Item(973, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(973, 9).Value  = 100104005
$ws.Cells.Item(973, 10).Value = "Pera"
$ws.Cells.Item(973, 11).Value = "Packham's Triumph"
$ws.Cells.Item(973, 12).Value = "Especial"
$ws.Cells.Item(973, 13).Value = 330
$ws.Cells.Item(973, 14).Value = 18000
$ws.Cells.Item(973, 15).Value = 18000
$ws.Cells.Item(973, 16).Value = 18000
$ws.Cells.Item(973, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(973, 18).Value = "Provincia de Linares"
$ws.Cells.Item(973, 19).Value = 1000
$ws.Cells.Item(973, 20).Value = 18

# --- New row 974 -----------------------------------------------------
$ws.Cells.Item(974, 1).Value  = 9
$ws.Cells.Item(974, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(974, 3).Value  = "Metropolitana"
$ws.Cells.Item(974, 4).Value  = 45013
$ws.Cells.Item(974, 5).Value  = 13
$ws.Cells.Item(974, 6).Value  = "Fruta"
$ws.Cells.Item(974, 7).Value  = 100104
$ws.Cells.Item(974, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(974, 9).Value  = 100104005
$ws.Cells.Item(974, 10).Value = "Pera"
$ws.Cells.Item(974, 11).Value = "Packham's Triumph"
$ws.Cells.Item(974, 12).Value = "Primera"
$ws.Cells.Item(974, 13).Value = 280
$ws.Cells.Item(974, 14).Value = 16000
$ws.Cells.Item(974, 15).Value = 16000
$ws.Cells.Item(974, 16).Value = 16000
$ws.Cells.Item(974, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(974, 18).Value = "Provincia de Linares"
$ws.Cells.Item(974, 19).Value = 889
$ws.Cells.Item(974, 20).Value = 18

# Make sure the date cells keep the same date/time number format the
# rest of column D uses (style copied down automatically by Insert, but
# set explicitly too so the serial values 45013 render/compare the same
# way as their neighbours).
$ws.Range("D973:D974").NumberFormat = $ws.Range("D975").NumberFormat
